# Se agrego menu popup a la tabla
# Elimina los renglones de infracciones que ya no aplican (limpieza de datos
# relacionada con el nuevo menu emergente de agregar/modificar/eliminar
# registros en la ventana de consulta de Articulos).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pruebita")

# Filas (numero de fila de Excel, 1 = encabezado) que deben eliminarse por
# completo. Se procesan de mayor a menor para que no se recorran los indices
# de las filas restantes durante el borrado.
$rowsToDelete = @(45, 19, 18, 16, 14, 13, 11, 8)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
